# Cover letter edit (16 January 2018):
#   "as an Embedded Engineer" -> "as an Hardware Suport Engineer"
# "Embedded" was selected and retyped as "Hardware Suport" (misspelling of
# "Support" kept verbatim, as in the source revision), which splits the
# original run into three runs and drops Word's "_GoBack" last-edit
# bookmark right after the freshly typed text. The pre-existing
# "_GoBack" bookmark elsewhere in the letter is removed (Word only ever
# keeps one), and the other existing bookmark ("_Hlk477116554") keeps its
# position but is renumbered as a side effect of the new bookmark being
# inserted earlier in the document.

$d = $word.ActiveDocument

# Locate the run to edit and the "Embedded" word inside it.
$full = $d.Content
$full.Find.Execute("as an Embedded Engineer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fullStart = $full.Start
$fullEnd = $full.End

$mid = $d.Content
$mid.Find.Execute("Embedded", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$midStart = $mid.Start
$midEnd = $mid.End

# Drop temporary bookmarks at every split point so later text edits can't
# re-merge the pieces, and so the positions stay valid as lengths change.
$d.Bookmarks.Add("ZEdit0", $d.Range($fullStart, $fullStart))
$d.Bookmarks.Add("ZEdit1", $d.Range($midStart, $midStart))
$d.Bookmarks.Add("ZEdit2", $d.Range($midEnd, $midEnd))
$d.Bookmarks.Add("ZEdit3", $d.Range($fullEnd, $fullEnd))

function Set-Segment($markA, $markB, $text) {
    $a = $d.Bookmarks.Item($markA).Start
    $b = $d.Bookmarks.Item($markB).Start
    $r = $d.Range($a, $b)
    $r.Text = $text
}

# "Embedded" -> "Hardware Suport" (typed replacement becomes its own run).
Set-Segment "ZEdit1" "ZEdit2" "Hardware Suport"

# Re-type the trailing " Engineer" too so it becomes a fresh run of its
# own (a no-op Text= assignment back to identical content does not split
# a run off, so go through a placeholder first).
Set-Segment "ZEdit2" "ZEdit3" "ZZPLACEHOLDERZZ"
Set-Segment "ZEdit2" "ZEdit3" " Engineer"

# The leading "as an " keeps the identity of the original run untouched.

# Word keeps only one "_GoBack" bookmark - remove the old one...
$d.Bookmarks.Item("_GoBack").Delete()

# ...and drop the new one right after the newly typed text.
$goBackPos = $d.Bookmarks.Item("ZEdit2").Start
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

# Clean up the scaffolding bookmarks used to track positions.
$d.Bookmarks.Item("ZEdit0").Delete()
$d.Bookmarks.Item("ZEdit1").Delete()
$d.Bookmarks.Item("ZEdit2").Delete()
$d.Bookmarks.Item("ZEdit3").Delete()
